$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 78.67
$ws.Range("G3").Value = 76.75
$ws.Range("G4").Value = 77.12
$ws.Range("G5").Value = 75.06
$ws.Range("G6").Value = 81.11
$ws.Range("G7").Value = 81.78
$ws.Range("G8").Value = 81.90000000000001
$ws.Range("G9").Value = 80.97
$ws.Range("G10").Value = 80.97
$ws.Range("G12").Value = 84.65000000000001
$ws.Range("G13").Value = 80.83
$ws.Range("G14").Value = 83.38
$ws.Range("G15").Value = 84.17
$ws.Range("G16").Value = 84.56
$ws.Range("G17").Value = 76.87
$ws.Range("G18").Value = 84.73999999999999
$ws.Range("G19").Value = 84.92
$ws.Range("G20").Value = 83.03
$ws.Range("G21").Value = 83.81
$ws.Range("G22").Value = 81.56
$ws.Range("G25").Value = 80.37
$ws.Range("G26").Value = 82.72
$ws.Range("G27").Value = 82.72
$ws.Range("G28").Value = 82.72
$ws.Range("G29").Value = 82.72
$ws.Range("G30").Value = 82.72
$ws.Range("G33").Value = 82.72
$ws.Range("G35").Value = 83.36
$ws.Range("G36").Value = 81.94
$ws.Range("G38").Value = 80.97
$ws.Range("G39").Value = 78.98999999999999
$ws.Range("G40").Value = 73.25
$ws.Range("G41").Value = 83
$ws.Range("G42").Value = 84.89
$ws.Range("G43").Value = 83.56999999999999
$ws.Range("G44").Value = 75.48999999999999
$ws.Range("G45").Value = 81.19
$ws.Range("G46").Value = 77.88
$ws.Range("G47").Value = 74.89
$ws.Range("G48").Value = 73.02
$ws.Range("G49").Value = 72.66
$ws.Range("G50").Value = 69.78
$ws.Range("G51").Value = 78.61
$ws.Range("G52").Value = 70.28
$ws.Range("G53").Value = 73.09999999999999
$ws.Range("G54").Value = 78.72
$ws.Range("G55").Value = 73.54000000000001
$ws.Range("G56").Value = 63.58
$ws.Range("G57").Value = 77.09999999999999
$ws.Range("G58").Value = 63.44
$ws.Range("G59").Value = 81.98
$ws.Range("G60").Value = 77.05
$ws.Range("G61").Value = 67.36
$ws.Range("G62").Value = 70.05
$ws.Range("G63").Value = 79.54000000000001
$ws.Range("G64").Value = 61.59
$ws.Range("G65").Value = 65.59
$ws.Range("G66").Value = 76.44
$ws.Range("G67").Value = 86.94
$ws.Range("G68").Value = 77.16
$ws.Range("G69").Value = 76.14
$ws.Range("G70").Value = 71.37
$ws.Range("G71").Value = 75.59999999999999
$ws.Range("G72").Value = 74.62
$ws.Range("G73").Value = 78.56
$ws.Range("G74").Value = 68.36
$ws.Range("G75").Value = 72.27
$ws.Range("G76").Value = 82.39
$ws.Range("G77").Value = 71.16
$ws.Range("G78").Value = 92.73
$ws.Range("G79").Value = 68.90000000000001
$ws.Range("G80").Value = 89.55
$ws.Range("G81").Value = 72.5
$ws.Range("G82").Value = 73.51000000000001
$ws.Range("G83").Value = 67.44
